$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "navprice" column (C) to hold
# the new "fnd_ver" field pulled from the database comparison.
$ws.Range("C1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("C1").Value = "fnd_ver"

# Fill the new column with the fund version value for every data row.
# The leading apostrophe forces it to be stored as text (matching the
# quote-prefixed text style already used by the neighboring fnd_id column),
# rather than being interpreted as the number 1.
$ws.Range("C2:C6").Value = "'01"

# Reflect the new column in the current selection.
$ws.Range("C3:C6").Select()
